# Prepend "Design: " to the start of each answer paragraph's text in the
# feedback table. Each target string below uniquely identifies the start of
# one of the six answer paragraphs.

$d = $word.ActiveDocument

$targets = @(
    "Way too many systems",
    "More-less straight forward",
    "I consider material handling",
    "Block manufacturing timetable",
    "In my opinion communication",
    "We got us well employed"
)

foreach ($t in $targets) {
    $d.Content.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "Design: $t", 2) | Out-Null
}
